# Vuvlo_calc.xlsx edit: update 1.25k resistors to 1.2k resistors
#
# - "9V" sheet: Rfbb (B6) changes from 9 to 9.33, which recalculates the
#   dependent Rfbt formula in B4 (B3/(B6-1)) from 1.25 to ~1.2004801920768307.
# - "PowerConditioner" sheet: two new rows are appended documenting the new
#   SETI resistor value (2.4kOhm).
# - Selection / active-sheet bookkeeping is updated to match where the
#   author was last working (PowerConditioner sheet, cell A20) instead of
#   the previous active sheet ("9V", cell B11).

$wb = $excel.ActiveWorkbook

$wsPower = $wb.Worksheets.Item("PowerConditioner")
$ws9V = $wb.Worksheets.Item("9V")

# --- 9V sheet: Rfbb 9 -> 9.33 (ripples into Rfbt formula result) ---
$ws9V.Range("B6").Value = 9.33

# --- PowerConditioner sheet: document new SETI resistor ---
$wsPower.Range("A18").Value = "SETI"
$wsPower.Range("A19").Value = "2.4kOhm"

# --- view bookkeeping: move the active tab / selection ---
$ws9V.Range("B7").Select()

$wsPower.Activate()
$wsPower.Range("A20").Select()
